$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# The "Status" column on the Overview sheet shares the same string as the
# per-language "Status" cell (both were "Ready for handoff"). Update it here
# too so the old text is fully replaced everywhere, matching the new report.
$overview.Range("B2").Value = "Handoff transform failed"
$overview.Range("C2").Value = "Handoff transform failed"

foreach ($ws in @($zhcn, $dede)) {
    # Status -> "Handoff transform failed"
    $ws.Range("B2").Value = "Handoff transform failed"

    # Remove the "Latest Handoff File" hyperlink + its cell value entirely
    # (the report no longer has a handoff file to link to).
    $links = @()
    foreach ($h in $ws.Hyperlinks) { $links += $h }
    foreach ($h in $links) {
        if ($h.Range.Address() -eq "`$C`$2") {
            $h.Delete()
        }
    }
    $ws.Range("C2").Clear()

    # Latest Handoff Datetime reset to the zero-date sentinel
    $ws.Range("D2").Value = "0001-01-01 00:00:00"

    # Handoff Reason -> "Ignored"
    $ws.Range("H2").Value = "Ignored"
}
